$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 1014; this shifts existing rows 1014:1067
# down to 1016:1069 and extends the sheet dimension accordingly.
$ws.Rows("1014:1015").Insert()

# Populate the newly inserted row 1014 (Primera quality, new weekly entry).
$ws.Range("A1014").Value = 3
$ws.Range("B1014").Value = "Femacal de La Calera"
$ws.Range("C1014").Value = "Coquimbo"
$ws.Range("D1014").Value = 45041
$ws.Range("E1014").Value = 5
$ws.Range("F1014").Value = 100114014
$ws.Range("G1014").Value = "Betarraga"
$ws.Range("H1014").Value = "Sin especificar"
$ws.Range("I1014").Value = "Primera"
$ws.Range("J1014").Value = 3700
$ws.Range("K1014").Value = 750
$ws.Range("L1014").Value = 800
$ws.Range("M1014").Value = 776
$ws.Range("N1014").Value = "$/paquete 4 unidades"
$ws.Range("O1014").Value = "Provincia de Quillota"
$ws.Range("P1014").Value = 194
$ws.Range("Q1014").Value = 4
$ws.Range("R1014").Value = "Hortaliza"

# Populate the newly inserted row 1015 (Segunda quality, new weekly entry).
$ws.Range("A1015").Value = 3
$ws.Range("B1015").Value = "Femacal de La Calera"
$ws.Range("C1015").Value = "Coquimbo"
$ws.Range("D1015").Value = 45041
$ws.Range("E1015").Value = 5
$ws.Range("F1015").Value = 100114014
$ws.Range("G1015").Value = "Betarraga"
$ws.Range("H1015").Value = "Sin especificar"
$ws.Range("I1015").Value = "Segunda"
$ws.Range("J1015").Value = 1950
$ws.Range("K1015").Value = 600
$ws.Range("L1015").Value = 600
$ws.Range("M1015").Value = 600
$ws.Range("N1015").Value = "$/paquete 4 unidades"
$ws.Range("O1015").Value = "Provincia de Quillota"
$ws.Range("P1015").Value = 150
$ws.Range("Q1015").Value = 4
$ws.Range("R1015").Value = "Hortaliza"
